$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helpers: write a cell's value while forcing the exact same cell style
# (cellXf index) as an existing, style-stable "donor" cell elsewhere on the
# sheet. A plain ".Value = x" assignment never changes the number format /
# font / border of the target cell, but on its own it also can't switch a
# cell between "number" and "text-that-looks-like-a-number" (e.g. the "0"
# and "***.*" placeholders that are stored as shared strings, not numbers).
# Copy + PasteSpecial(xlPasteFormats) reapplies the donor's exact existing
# cellXf instead of minting a brand-new style, which is what a bare
# NumberFormat assignment would do.
# ---------------------------------------------------------------------------

function Set-NumCell($cellRef, $styleSrcRef, $val) {
    $dst = $ws.Range($cellRef)
    $src = $ws.Range($styleSrcRef)
    $dst.Value = $val
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

function Set-TextCell($cellRef, $styleSrcRef, $text) {
    $dst = $ws.Range($cellRef)
    $src = $ws.Range($styleSrcRef)
    $dst.NumberFormat = "@"
    $dst.Value = $text
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Header: new report volume number / week-covering date range
# ---------------------------------------------------------------------------

$ws.Range("A8").Value = "Volume 32   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/7/2025  Through  7/13/2025"

# ---------------------------------------------------------------------------
# Crime-complaint stat grid updates (rows 15-31)
# ---------------------------------------------------------------------------

Set-NumCell "N15" "N14" -80

Set-NumCell "G16" "G15" 4
Set-NumCell "J16" "G15" 25
Set-NumCell "K16" "N14" -92
Set-NumCell "M16" "N14" -86.666666666666
Set-NumCell "N16" "N14" -97.849462365591

Set-TextCell "D17" "C15" "0"
Set-TextCell "E17" "E15" "***.*"
Set-NumCell "G17" "G15" 1
Set-NumCell "H17" "N14" 100
Set-NumCell "I17" "G15" 9
Set-NumCell "K17" "N14" 12.5
Set-NumCell "L17" "N14" 80
Set-NumCell "M17" "N14" 800
Set-NumCell "N17" "N14" -59.090909090909

Set-NumCell "C19" "G15" 2
Set-NumCell "D19" "G15" 1
Set-NumCell "E19" "N14" 100
Set-NumCell "F19" "G15" 6
Set-NumCell "G19" "G15" 4
Set-NumCell "H19" "N14" 50
Set-NumCell "I19" "G15" 23
Set-NumCell "J19" "G15" 28
Set-NumCell "K19" "N14" -17.857142857142
Set-NumCell "L19" "N14" 0
Set-NumCell "M19" "N14" -36.111111111111
Set-NumCell "N19" "N14" -74.444444444444

Set-NumCell "C21" "D21" 2
Set-NumCell "E21" "M21" 0
Set-NumCell "F21" "D21" 8
Set-NumCell "G21" "D21" 10
Set-NumCell "H21" "M21" -20
Set-NumCell "I21" "D21" 35
Set-NumCell "J21" "D21" 63
Set-NumCell "K21" "M21" -44.444444444444
Set-NumCell "L21" "M21" -7.894736842105
Set-NumCell "M21" "M21" -37.5
Set-NumCell "N21" "M21" -84.978540772532

Set-TextCell "C24" "C15" "0"
Set-TextCell "D24" "C15" "0"
Set-TextCell "E24" "E15" "***.*"
Set-NumCell "F24" "G15" 4
Set-NumCell "G24" "G15" 3
Set-NumCell "H24" "N14" 33.333333333333
Set-NumCell "L24" "N14" -15.789473684210
Set-NumCell "M24" "N14" -56.756756756756

Set-TextCell "F25" "C15" "0"

Set-NumCell "C26" "G15" 1
Set-TextCell "D26" "C15" "0"
Set-TextCell "E26" "E15" "***.*"
Set-NumCell "F26" "G15" 4
Set-NumCell "G26" "G15" 8
Set-NumCell "L26" "N14" -44.827586206896

Set-NumCell "C28" "G15" 1
Set-NumCell "F28" "G15" 5
Set-NumCell "I28" "G15" 13
Set-NumCell "K28" "N14" 333.333333333333
Set-NumCell "L28" "N14" 8.333333333333

Set-NumCell "D31" "G15" 1
Set-NumCell "E31" "N14" -100
Set-NumCell "G31" "G15" 1
Set-NumCell "H31" "N14" -100
Set-NumCell "J31" "G15" 1
Set-NumCell "K31" "N14" -100
